$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = -1
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 1594.12
$ws.Cells.Item(2, 7).Value = 19.65
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 25).Value = [double]"1.034588057504798e-16"
$ws.Cells.Item(2, 26).Value = [double]"2.433896954542872e-12"
# Row 3
$ws.Cells.Item(3, 4).Value = 999817
$ws.Cells.Item(3, 5).Value = 1640
$ws.Cells.Item(3, 6).Value = 1593.2
$ws.Cells.Item(3, 7).Value = 20.04
$ws.Cells.Item(3, 11).Value = [double]"2.75e-05"
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 16).Value = [double]"9.999208758039257e-07"
$ws.Cells.Item(3, 17).Value = [double]"4.516455053091704e-11"
$ws.Cells.Item(3, 18).Value = 2.023240151271259
$ws.Cells.Item(3, 19).Value = [double]"1.629936922501878e-06"
$ws.Cells.Item(3, 20).Value = [double]"3.789946124422143e-07"
$ws.Cells.Item(3, 21).Value = 2.218414459480953
$ws.Cells.Item(3, 22).Value = 0.1630116514466714
$ws.Cells.Item(3, 23).Value = 0.03790363750743998
$ws.Cells.Item(3, 24).Value = 2.2184144569907
$ws.Cells.Item(3, 25).Value = [double]"2.199836279894013e-11"
$ws.Cells.Item(3, 26).Value = [double]"2.702610592513617e-07"
# Row 4
$ws.Cells.Item(4, 4).Value = -1
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 1592.78
$ws.Cells.Item(4, 7).Value = 20.08
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 25).Value = [double]"1.233663930019391e-16"
$ws.Cells.Item(4, 26).Value = [double]"2.67506395485076e-12"
# Row 5
$ws.Cells.Item(5, 4).Value = 999944
$ws.Cells.Item(5, 5).Value = 396
$ws.Cells.Item(5, 6).Value = 1591.89
$ws.Cells.Item(5, 7).Value = 20.11
$ws.Cells.Item(5, 11).Value = [double]"2.75e-05"
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 16).Value = [double]"1.000047888871531e-07"
$ws.Cells.Item(5, 17).Value = [double]"4.552348543349149e-12"
$ws.Cells.Item(5, 18).Value = 2.023004113939691
$ws.Cells.Item(5, 19).Value = [double]"3.959846756507628e-08"
$ws.Cells.Item(5, 20).Value = [double]"3.784252074943758e-08"
$ws.Cells.Item(5, 21).Value = 2.218379278774329
$ws.Cells.Item(5, 22).Value = 0.03959780115212037
$ws.Cells.Item(5, 23).Value = 0.0378418839101956
$ws.Cells.Item(5, 24).Value = 2.218379278602178
$ws.Cells.Item(5, 25).Value = [double]"2.200228792112639e-12"
$ws.Cells.Item(5, 26).Value = [double]"2.700974446376493e-08"
# Row 6
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 1592.17
$ws.Cells.Item(6, 7).Value = 20.47
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 25).Value = [double]"1.25394704433805e-16"
$ws.Cells.Item(6, 26).Value = [double]"2.698351439949472e-12"
# Row 7
$ws.Cells.Item(7, 4).Value = 1000003
$ws.Cells.Item(7, 5).Value = 479
$ws.Cells.Item(7, 6).Value = 1592.19
$ws.Cells.Item(7, 7).Value = 20.21
$ws.Cells.Item(7, 11).Value = [double]"2.75e-05"
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 16).Value = [double]"1.000106894942308e-08"
$ws.Cells.Item(7, 17).Value = [double]"4.659875191787677e-13"
$ws.Cells.Item(7, 18).Value = 2.022888385509832
$ws.Cells.Item(7, 19).Value = [double]"4.749816185204857e-09"
$ws.Cells.Item(7, 20).Value = [double]"3.785890133436816e-09"
$ws.Cells.Item(7, 21).Value = 2.218386761343267
$ws.Cells.Item(7, 22).Value = 0.04749456014900783
$ws.Cells.Item(7, 23).Value = 0.03785603059362858
$ws.Cells.Item(7, 24).Value = 2.218386760987068
$ws.Cells.Item(7, 25).Value = [double]"2.201489100138555e-13"
$ws.Cells.Item(7, 26).Value = [double]"2.704073817686499e-09"
# Row 8
$ws.Cells.Item(8, 5).Value = 12
$ws.Cells.Item(8, 6).Value = 1592.13
$ws.Cells.Item(8, 7).Value = 20.61
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 25).Value = [double]"1.256178829618184e-16"
# Row 9
$ws.Cells.Item(9, 4).Value = 999707
$ws.Cells.Item(9, 5).Value = -218
$ws.Cells.Item(9, 6).Value = 1592.07
$ws.Cells.Item(9, 7).Value = 20.19
$ws.Cells.Item(9, 11).Value = [double]"9e-06"
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 16).Value = [double]"9.998728528891004e-10"
$ws.Cells.Item(9, 17).Value = [double]"3.565722513400956e-14"
$ws.Cells.Item(9, 18).Value = 2.076177969564305
$ws.Cells.Item(9, 19).Value = [double]"-2.300055890449676e-10"
$ws.Cells.Item(9, 20).Value = [double]"3.784447589842954e-10"
$ws.Cells.Item(9, 21).Value = 2.21837789659931
$ws.Cells.Item(9, 22).Value = -0.02300419819965971
$ws.Cells.Item(9, 23).Value = 0.03785046389708496
$ws.Cells.Item(9, 24).Value = 2.218377896754195
$ws.Cells.Item(9, 25).Value = [double]"2.212282538065619e-14"
$ws.Cells.Item(9, 26).Value = [double]"2.727548546400607e-10"
# Row 10
$ws.Cells.Item(10, 4).Value = 999876
$ws.Cells.Item(10, 5).Value = -158
$ws.Cells.Item(10, 6).Value = 1591.9
$ws.Cells.Item(10, 7).Value = 20.8
$ws.Cells.Item(10, 11).Value = [double]"1e-06"
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 16).Value = [double]"1.000041880926454e-09"
$ws.Cells.Item(10, 17).Value = [double]"3.428163893434152e-14"
$ws.Cells.Item(10, 18).Value = 2.10362861939283
$ws.Cells.Item(10, 19).Value = [double]"-1.700041310332369e-10"
$ws.Cells.Item(10, 20).Value = [double]"3.784281038733997e-10"
$ws.Cells.Item(10, 21).Value = 2.218376634738403
$ws.Cells.Item(10, 22).Value = -0.01700022913772975
$ws.Cells.Item(10, 23).Value = 0.03784240088320211
$ws.Cells.Item(10, 24).Value = 2.218376634837701
$ws.Cells.Item(10, 25).Value = [double]"2.212655741279608e-14"
$ws.Cells.Item(10, 26).Value = [double]"2.727713777108246e-10"
# Row 11
$ws.Cells.Item(11, 4).Value = 999859
$ws.Cells.Item(11, 5).Value = -138
$ws.Cells.Item(11, 6).Value = 1591.86
$ws.Cells.Item(11, 7).Value = 20.8
$ws.Cells.Item(11, 11).Value = [double]"1e-06"
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 16).Value = [double]"1.000024878106128e-09"
$ws.Cells.Item(11, 17).Value = [double]"3.428106579316784e-14"
$ws.Cells.Item(11, 18).Value = 2.103628376104718
$ws.Cells.Item(11, 19).Value = [double]"-1.500036450293267e-10"
$ws.Cells.Item(11, 20).Value = [double]"3.784028682873857e-10"
$ws.Cells.Item(11, 21).Value = 2.21837526293537
$ws.Cells.Item(11, 22).Value = -0.01500045722006737
$ws.Cells.Item(11, 23).Value = 0.03784052072088269
$ws.Cells.Item(11, 24).Value = 2.218375263012689
$ws.Cells.Item(11, 25).Value = [double]"2.212618650771077e-14"
$ws.Cells.Item(11, 26).Value = [double]"2.727599320707821e-10"
# Row 12
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = 30
$ws.Cells.Item(12, 6).Value = 1595.3
$ws.Cells.Item(12, 7).Value = 19.9
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 25).Value = [double]"1.253792702963777e-16"
$ws.Cells.Item(12, 26).Value = [double]"2.706498021048579e-12"
# Row 13
$ws.Cells.Item(13, 4).Value = 999917
$ws.Cells.Item(13, 5).Value = 2270
$ws.Cells.Item(13, 6).Value = 1595.12
$ws.Cells.Item(13, 7).Value = 20.02
$ws.Cells.Item(13, 11).Value = [double]"9e-06"
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 16).Value = [double]"1.000003386072989e-10"
$ws.Cells.Item(13, 17).Value = [double]"3.55239852036184e-15"
$ws.Cells.Item(13, 18).Value = 2.078393745352665
$ws.Cells.Item(13, 19).Value = [double]"2.239887559529839e-10"
$ws.Cells.Item(13, 20).Value = [double]"3.799313632024773e-11"
$ws.Cells.Item(13, 21).Value = 2.218469092673973
$ws.Cells.Item(13, 22).Value = 0.2239949506835178
$ws.Cells.Item(13, 23).Value = 0.03799418746939371
$ws.Cells.Item(13, 24).Value = 2.218469107537312
$ws.Cells.Item(13, 25).Value = [double]"2.325400866735839e-15"
$ws.Cells.Item(13, 26).Value = [double]"2.976685754945776e-11"
# Row 14
$ws.Cells.Item(14, 4).Value = 99984
$ws.Cells.Item(14, 5).Value = 240
$ws.Cells.Item(14, 6).Value = 1592.93
$ws.Cells.Item(14, 7).Value = 20
$ws.Cells.Item(14, 11).Value = [double]"9e-06"
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 16).Value = [double]"9.998813749229067e-12"
$ws.Cells.Item(14, 17).Value = [double]"5.530009498881907e-16"
$ws.Cells.Item(14, 18).Value = 2.000815441738256
$ws.Cells.Item(14, 19).Value = [double]"2.099894587059224e-11"
$ws.Cells.Item(14, 20).Value = [double]"5.577920745996919e-12"
$ws.Cells.Item(14, 21).Value = 2.018848486255812
$ws.Cells.Item(14, 22).Value = 0.2100208914940808
$ws.Cells.Item(14, 23).Value = 0.05578755672234557
$ws.Cells.Item(14, 24).Value = 2.01884848617405
$ws.Cells.Item(14, 25).Value = [double]"3.455396987398924e-16"
$ws.Cells.Item(14, 26).Value = [double]"5.404498724452729e-12"
# Row 15
$ws.Cells.Item(15, 4).Value = 10003
$ws.Cells.Item(15, 5).Value = 50
$ws.Cells.Item(15, 6).Value = 1594.62
$ws.Cells.Item(15, 7).Value = 19.96
$ws.Cells.Item(15, 11).Value = [double]"1.2e-05"
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 17).Value = [double]"4.427353131389402e-16"
$ws.Cells.Item(15, 18).Value = 2.021871798193535
$ws.Cells.Item(15, 19).Value = [double]"1.999899606723071e-12"
$ws.Cells.Item(15, 20).Value = [double]"4.425964763871571e-12"
$ws.Cells.Item(15, 21).Value = 2.021927365067199
$ws.Cells.Item(15, 22).Value = 0.2000178963354084
$ws.Cells.Item(15, 23).Value = 0.4426582965782722
$ws.Cells.Item(15, 24).Value = 2.02192735994078
$ws.Cells.Item(15, 25).Value = [double]"1.474303418427843e-16"
$ws.Cells.Item(15, 26).Value = [double]"2.975835895206265e-12"
# Row 16
$ws.Cells.Item(16, 4).Value = 99992
$ws.Cells.Item(16, 5).Value = 225
$ws.Cells.Item(16, 6).Value = 1594.5
$ws.Cells.Item(16, 7).Value = 20.6
$ws.Cells.Item(16, 11).Value = [double]"1e-06"
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 16).Value = [double]"9.999613822344359e-12"
$ws.Cells.Item(16, 17).Value = [double]"5.439387857512342e-16"
$ws.Cells.Item(16, 18).Value = 2.003742609396941
$ws.Cells.Item(16, 19).Value = [double]"1.949902116554994e-11"
$ws.Cells.Item(16, 20).Value = [double]"5.581826060699703e-12"
$ws.Cells.Item(16, 21).Value = 2.01896089660569
$ws.Cells.Item(16, 22).Value = 0.1950037960594554
$ws.Cells.Item(16, 23).Value = 0.05582214911600104
$ws.Cells.Item(16, 24).Value = 2.018960896683757
$ws.Cells.Item(16, 25).Value = [double]"3.454336747425682e-16"
$ws.Cells.Item(16, 26).Value = [double]"5.410041846408926e-12"
# Row 17
$ws.Cells.Item(17, 4).Value = 999873
$ws.Cells.Item(17, 5).Value = 2265
$ws.Cells.Item(17, 6).Value = 1593.06
$ws.Cells.Item(17, 7).Value = 20.8
$ws.Cells.Item(17, 11).Value = [double]"1e-06"
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 16).Value = [double]"9.999593820516475e-11"
$ws.Cells.Item(17, 17).Value = [double]"3.427856746129195e-15"
$ws.Cells.Item(17, 18).Value = 2.103634705781022
$ws.Cells.Item(17, 19).Value = [double]"2.234887810513032e-10"
$ws.Cells.Item(17, 20).Value = [double]"3.789425685746423e-11"
$ws.Cells.Item(17, 21).Value = 2.218414622504492
$ws.Cells.Item(17, 22).Value = 0.2235047970472111
$ws.Cells.Item(17, 23).Value = 0.03789697285417683
$ws.Cells.Item(17, 24).Value = 2.218414639620134
$ws.Cells.Item(17, 25).Value = [double]"2.325466191002344e-15"
$ws.Cells.Item(17, 26).Value = [double]"2.972722623579457e-11"
